$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.767793560937871
$ws.Range("C2").Value = 7.51257437759147
$ws.Range("D2").Value = 7.931057780663799

$ws.Range("B3").Value = 1.383656700623011
$ws.Range("C3").Value = 1.955583791664621
$ws.Range("D3").Value = 2.393202307189311

$ws.Range("B4").Value = 0.1190132055564578
$ws.Range("C4").Value = 0.1636774942882221
$ws.Range("D4").Value = 0.1969453475959414

$ws.Range("B5").Value = 89.74005669753598
$ws.Range("C5").Value = 90.67258625751462
$ws.Range("D5").Value = 91.14213403237282
